$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Subtract a 2-hour timedelta from every time value in column F (rows 2-29)
$timedelta = 2 / 24

for ($r = 2; $r -le 29; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $cell.Value2 = $cell.Value2 - $timedelta
}
